# Fix data type conversion for values with nested quotes
#
# The shared string "dolor" (used by cell E4) needs to become a longer
# string that contains embedded double quotes, used to test escaping of
# quoted values during data-type conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in E4 ("dolor" -> "dolor with "quoted" string to test escaping")
$ws.Range("E4").Value = 'dolor with "quoted" string to test escaping'

# Widen column E to fit the new, longer text
$ws.Columns("E").ColumnWidth = 40.8

# Update the active selection to E5 (as left by the author after editing E4)
[void]$ws.Range("E5").Select()
